$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected (legacy password hash) and the diff does not touch
# the <sheetProtection> element itself, so we must NOT call Protect/Unprotect on
# the sheet (that would drop/rehash the original protection). Instead, for each
# cell we need to edit, we temporarily unlock just that cell, write the new
# value/text, then relock it - the sheet-level protection XML stays untouched.

function Set-ProtectedValue {
    param($addr, $val)
    $c = $ws.Range($addr)
    $c.Locked = $false
    $c.Value = $val
    $c.Locked = $true
}

# 1) Update the confidentiality / "as of" date footer text (A18) from
#    2021-05-17 to 2021-05-18.
$oldText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-17 for illustrative purposes only and are subject to change."
$newText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-18 for illustrative purposes only and are subject to change."

$a18 = $ws.Range("A18")
if ($a18.Value -eq $oldText) {
    Set-ProtectedValue "A18" $newText
}

# 2) Update the Weight (D) and Percent Change (E) values for rows 2-15.
Set-ProtectedValue "D2"  0.05650056643735531
Set-ProtectedValue "E2"  -0.007811047052259656

Set-ProtectedValue "D3"  0.02364325939547178
Set-ProtectedValue "E3"  0.003499902780478292

Set-ProtectedValue "D4"  0.03063490239217167
Set-ProtectedValue "E4"  0.01447038394752065

Set-ProtectedValue "D5"  0.03379053310785295
Set-ProtectedValue "E5"  -0.01477104874446078

Set-ProtectedValue "D6"  0.03830101630354475
Set-ProtectedValue "E6"  0.004732510288065761

Set-ProtectedValue "D7"  0.01932191702103207
Set-ProtectedValue "E7"  -0.008327385201047033

Set-ProtectedValue "D8"  0.00421626414041079
Set-ProtectedValue "E8"  0.01869158878504695

Set-ProtectedValue "D9"  0.006890505196447043
Set-ProtectedValue "E9"  -0.01372474266107504

Set-ProtectedValue "D10" 0.07305564220861312
Set-ProtectedValue "E10" 0.001618122977346426

Set-ProtectedValue "D11" 0.07317385522189568
Set-ProtectedValue "E11" 0.001077005923532548

Set-ProtectedValue "D12" 0.1437365163281724
Set-ProtectedValue "E12" -0.002558666569193568

Set-ProtectedValue "D13" 0.3819620076510089
Set-ProtectedValue "E13" 0

Set-ProtectedValue "D14" 0.1147730145960234
Set-ProtectedValue "E14" -0.004085556356644982

Set-ProtectedValue "D15" 0.9999999999999999
Set-ProtectedValue "E15" -0.001049468862364011
